$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Date line
Replace-Text "2025-02-22 Saturday" "2025-02-23 Sunday"

# Problem table cells - the value "34÷4=" repeats twice, so handle the table
# directly by row/column position to avoid ambiguity; all other values are
# unique in the document and are replaced via Find/Replace.
$table = $d.Tables.Item(1)

# Row 1 (table row index 1)
Replace-Text "97÷3=" "38÷8="
Replace-Text "18÷3=" "12÷3="
Replace-Text "45÷9=" "80÷3="
Replace-Text "53÷6=" "90÷5="
Replace-Text "43÷8=" "52÷6="

# Row 5 (table row index 5)
Replace-Text "58÷3=" "39÷7="
Replace-Text "14÷6=" "24÷8="
Replace-Text "45÷4=" "48÷4="
Replace-Text "76÷4=" "76÷2="
Replace-Text "77÷9=" "73÷8="

# Row 9 (table row index 9) - last cell "34÷4=" is the first occurrence
Replace-Text "38÷6=" "10÷5="
Replace-Text "64÷3=" "28÷5="
Replace-Text "54÷7=" "45÷5="
Replace-Text "34÷5=" "45÷7="
$table.Cell(9, 5).Range.Text = "88÷5="

# Row 13 (table row index 13)
Replace-Text "53÷3=" "61÷7="
Replace-Text "21÷7=" "92÷3="
Replace-Text "72÷4=" "81÷3="
Replace-Text "57÷8=" "51÷3="
Replace-Text "16÷9=" "48÷6="

# Row 17 (table row index 17) - third cell "34÷4=" is the second occurrence
Replace-Text "79÷9=" "80÷2="
Replace-Text "94÷5=" "35÷8="
$table.Cell(17, 3).Range.Text = "89÷4="
Replace-Text "31÷3=" "18÷2="
Replace-Text "66÷3=" "91÷2="
